# Marksheet update: re-mark the quiz with a new question set / scoring
# (handles float-valued marking scheme inputs without corrupting adjacent
# cells - e.g. the old "-1" wrong-answer mark had been written as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 10-12) ------------------------------------------

$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 10
$ws.Range("E10").Value = 28

$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "62/112"

# --- Drop the third Student/Correct Ans pair (columns G:H) ---------------

$ws.Range("G15:H40").Clear()

# --- Drop the per-question detail for rows 19-40 in columns D:E ----------
# (only rows 16-18 keep a second Student/Correct Ans pair going forward)

$ws.Range("D19:E40").Clear()

# --- Second Student/Correct Ans pair, rows 16-18 --------------------------

$ws.Range("D16").Style = "correctStyle"
$ws.Range("D16").Value = "Option A"

$ws.Range("D17").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"

$ws.Range("D18").Style = "incorrectStyle"
$ws.Range("D18").Value = "Option B"

# --- Student answers, column A (correct -> green, wrong -> red) ----------

$ws.Range("A17").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"

$ws.Range("A18").Style = "incorrectStyle"
$ws.Range("A18").Value = "Option A"

$ws.Range("A19").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"

$ws.Range("A21").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"

$ws.Range("A22").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"

$ws.Range("A23").Style = "correctStyle"
$ws.Range("A23").Value = "Option D"

$ws.Range("A24").Style = "correctStyle"
$ws.Range("A24").Value = "Option A"

$ws.Range("A25").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"

$ws.Range("A32").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"

$ws.Range("A33").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"

$ws.Range("A34").Style = "correctStyle"
$ws.Range("A34").Value = "Option B"

$ws.Range("A36").Style = "correctStyle"
$ws.Range("A36").Value = "Option A"

$ws.Range("A37").Style = "correctStyle"
$ws.Range("A37").Value = "Option A"

$ws.Range("A38").Style = "correctStyle"
$ws.Range("A38").Value = "Option A"

$ws.Range("A39").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
